$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 175, shifting rows 175:248 down to 176:249
$ws.Rows.Item(175).Insert()

# Populate the new row 175 with the new weekly data point
$ws.Cells.Item(175, 1).Value = 5
$ws.Cells.Item(175, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(175, 3).Value = "Maule"
$ws.Cells.Item(175, 4).Value = 44839
$ws.Cells.Item(175, 5).Value = 7
$ws.Cells.Item(175, 6).Value = 100112017
$ws.Cells.Item(175, 7).Value = "Apio"
$ws.Cells.Item(175, 8).Value = "Americana (o)"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 700
$ws.Cells.Item(175, 11).Value = 9000
$ws.Cells.Item(175, 12).Value = 9000
$ws.Cells.Item(175, 13).Value = 9000
$ws.Cells.Item(175, 14).Value = "$/docena de matas"
$ws.Cells.Item(175, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(175, 16).Value = 1500
$ws.Cells.Item(175, 17).Value = 6
$ws.Cells.Item(175, 18).Value = "Hortaliza"

# Match the date style used by the rest of column D
$ws.Cells.Item(175, 4).NumberFormat = $ws.Cells.Item(176, 4).NumberFormat
